$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 13 (Docentes responsaveis value moves here)
$ws.Rows.Item(13).Insert()

# Fix the new row: remove the leftover styled cell in column A, and set B/C
# with the correct (non-bold, wrapped) body style before assigning values.
$ws.Cells.Item(13, 1).Clear()
$ws.Cells.Item(13, 2).Font.Bold = $false
$ws.Cells.Item(13, 2).WrapText = $true
$ws.Cells.Item(13, 2).VerticalAlignment = -4160
$ws.Cells.Item(13, 3).WrapText = $true
$ws.Cells.Item(13, 3).VerticalAlignment = -4160
$ws.Cells.Item(13, 3).Font.Color = 255

$docente = '5701460 - Antonio Iacono'
$ws.Cells.Item(13, 2).Value = $docente
$ws.Cells.Item(13, 3).Value = $docente

# Row 10: Objetivos (PT) gets the real paragraph instead of the misplaced docente text
$objetivosPt = 'Apresentar aspectos relacionados à inovação tecnológica para a competitividade das empresas, explorando os principais elementos conceituais que orientam as estratégias de criação e gestão da inovação, e suas implicações para gestão organizacional.'
$ws.Cells.Item(10, 2).Value = $objetivosPt
$ws.Cells.Item(10, 3).Value = $objetivosPt

# Row 14: Programa resumido (PT short syllabus) replaces "Semestral"
$programaResumidoPt = '1. Conceitos e Definições no Estudo da Inovação. 2. Tipos de inovação. 3. Ambientes de inovação. 4. Sistemas de Inovação. 5. Capacidade tecnológica e aprendizagem. 6. Processos de inovação. 7. Estratégias de inovação. 8. Políticas públicas para promoção da inovação. 9. Propriedade Intelectual.'
$ws.Cells.Item(14, 2).Value = $programaResumidoPt
$ws.Cells.Item(14, 3).Value = $programaResumidoPt

# Row 16: Programa (PT full syllabus) replaces the misplaced date value
$programaPt = '. Conceitos e Definições no Estudo da Inovação. 2. Tipos de inovação. 3. Ambientes de inovação. 4. Sistemas de Inovação. 5. Capacidade tecnológica e aprendizagem. 6. Processos de inovação. 7. Estratégias de inovação. 8. Políticas públicas para promoção da inovação. 9. Propriedade Intelectual.'
$ws.Cells.Item(16, 2).Value = $programaPt
$ws.Cells.Item(16, 3).Value = $programaPt

# Row 19: Metodo -> "Provas e Trabalhos"
$provasTrabalhos = 'Provas e Trabalhos'
$ws.Cells.Item(19, 2).Value = $provasTrabalhos
$ws.Cells.Item(19, 3).Value = $provasTrabalhos

# Row 20: Criterio -> grading formula
$criterioFormula = 'M = (0,8P + 0,2T)P = média aritmética de duas provas escritasT = Média das notas de trabalhos e exercíciosM = Média de aproveitamento do alunoAprovação com média de aproveitamento maior ou igual a 5,0 e no mínimo 70% de frequência às aulas.'
$ws.Cells.Item(20, 2).Value = $criterioFormula
$ws.Cells.Item(20, 3).Value = $criterioFormula

# Row 21: Norma de recuperacao -> recovery formula
$normaFormula = 'MF = (0,5 M + 0,5 R)M = Média de aproveitamento do aluno, antes da recuperaçãoR = Nota de uma prova de recuperaçãoMF = nota final de aproveitamento, após a recuperaçãoAprovação com média final de aproveitamento maior ou igual a 5,0.A recuperação deverá consistir de uma prova escrita englobando a matéria toda do semestre.Terá direito à prova de recuperação aqueles alunos reprovados com nota acima de 3,0 e frequência mínima de 70%.'
$ws.Cells.Item(21, 2).Value = $normaFormula
$ws.Cells.Item(21, 3).Value = $normaFormula

# Row 22: Bibliografia -> full reading list
$bibliografia = 'BARNEY, J.B.; CLARK, D.N. Resource-Based Theory: Creating and Sustaining Competitive Advantage. Oxford University Press, 2007.BESSANT, J.; TIDD, J. Inovação e empreendedorismo. Porto Alegre, Bookman, 2009.BURGELMAN, R. A.; CHRISTENSEN, C. M.; WHEELWRIGTH, S. C. Gestão estratégica da tecnologia e da inovação: conceitos e soluções. AMGH Editora, 2013.CONWAY, S; STEWARD, F. Managing and shaping innovation. Oxford University Press, 2009.CHRISTENSEN, Clayton M. O dilema da inovação. São Paulo: Makron Books, 2011.DAVILA, T; EPSTEIN, M. J.; SHELTON, R. As regras da Inovação. Porto Alegre, Bookman, 2008.DE NEGRI, J.A; SALERNO, M.S., (Orgs.). Inovação, padrões tecnológicos e desempenho das firmas industriais brasileiras. Brasília, Ipea, 2005.DODGSON, M.; GANN, D.; SALTER, A. The management of technological innovation: strategy and practice. Oxford University Press, 2008.DRUCKER, P.F. Inovação e espírito empreendedor. São Paulo: Pioneira, 2000.FIGUEIREDO, P.N. Gestão da inovação: conceitos, métricas e experiências de empresas no Brasil. Rio de Janeiro, LTC, 2015.FITZGERALD, E. et al. Inside Real Innovation: How the Right Approach Can Move Ideas from R&D to Market-And Get the Economy Moving. World Scientific, 2011.GOFFIN, K.; MITCHELL, R. Innovation management. 2nd ed. Palgrave – MacMillan, Houndsmill, 2010.HELFAT, C.E. et al. Dynamic capabilities: understanding strategic change in organizations. Blackwell Publishing, 2007.PRAHALAD,C.K.; KRISHNAN,M.S. The new of innovation. EUA: Editora Soundview Executive Book Sumaries, 2008.PROENÇA, A. et al. Gestão da inovação e competitividade no Brasil: da teoria para a prática. Bookman Editora, 2015.SALERNO, M.S.; GOMES, L.A.V. Gestão da inovação (mais) radical. Rio de Janeiro: Elsevier, 2018.SCHILLING, M.A. Strategic management of technological innovation. MacGraw-Hill/Irwin, 2009. TEECE, D. Capabilities and strategic management. In: Edited by Foss. N. Resources firms and strategies. A reader in the Resource-based Perspective. WA: Ed. Oxford University, 1987.TIDD, J.; BESSANT, J. Gestão da Inovação. Porto Alegre, Bookman, 2015.TIDD, J.; BESSANT, J. Strategic innovation management, Wiley, 2014.TIGRE, P. B. Gestão da inovação. Rio de Janeiro, Campus-Elsevier, 2006.TROTT, P. innovation management and new product development. Prentice Hall, 2008.WHITE, M. A.; BRUTON, G.D. The management of technology and innovation: a strategic approach. South-Western, Cengage Learning, 2011.'
$ws.Cells.Item(22, 2).Value = $bibliografia
$ws.Cells.Item(22, 3).Value = $bibliografia

Write-Host "done"
